# Update row 8 (year 2025) metrics on the active sheet to reflect the
# refreshed "bibi" data (commit: "atualizei dados da bibi e add")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = 1419
$ws.Range("D8").Value = 216
$ws.Range("E8").Value = 1203
$ws.Range("F8").Value = 8.859721082854799
$ws.Range("G8").Value = 84.77801268498943
$ws.Range("H8").Value = 15.22198731501057
